$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 6450
$ws.Cells.Item(51, 9).Value = 3375
$ws.Cells.Item(51, 10).Value = 8500
$ws.Cells.Item(51, 11).Value = 3375
$ws.Cells.Item(51, 12).Value = 8500
$ws.Cells.Item(51, 13).Value = -2891
$ws.Cells.Item(51, 14).Value = -9468

$ws.Cells.Item(86, 8).Value = 5292323.5
$ws.Cells.Item(86, 9).Value = 1583.15
$ws.Cells.Item(86, 10).Value = 11170924
$ws.Cells.Item(86, 11).Value = 1583.15
$ws.Cells.Item(86, 12).Value = 11170924
$ws.Cells.Item(86, 13).Value = -460.1500000000001
$ws.Cells.Item(86, 14).Value = -11173170

$ws.Cells.Item(89, 8).Value = 5292323.5
$ws.Cells.Item(89, 9).Value = 1583.15
$ws.Cells.Item(89, 10).Value = 11170924
$ws.Cells.Item(89, 11).Value = 7915.75
$ws.Cells.Item(89, 12).Value = 55854620
$ws.Cells.Item(89, 13).Value = -2299.75
$ws.Cells.Item(89, 14).Value = -55865852

$ws.Cells.Item(92, 8).Value = 104859.84
$ws.Cells.Item(92, 9).Value = 214.3125
$ws.Cells.Item(92, 10).Value = 662969.3
$ws.Cells.Item(92, 11).Value = 214.3125
$ws.Cells.Item(92, 12).Value = 662969.3
$ws.Cells.Item(92, 13).Value = 1033.6875
$ws.Cells.Item(92, 14).Value = -665465.3

$ws.Cells.Item(100, 8).Value = 8255
$ws.Cells.Item(100, 9).Value = 10499
$ws.Cells.Item(100, 10).Value = 2645
$ws.Cells.Item(100, 11).Value = 10499
$ws.Cells.Item(100, 12).Value = 2645
$ws.Cells.Item(100, 13).Value = -9958
$ws.Cells.Item(100, 14).Value = -3727

$ws.Cells.Item(118, 8).Value = 304.2857
$ws.Cells.Item(118, 9).Value = 255
$ws.Cells.Item(118, 10).Value = 600
$ws.Cells.Item(118, 11).Value = 765
$ws.Cells.Item(118, 12).Value = 1800
$ws.Cells.Item(118, 13).Value = 892
$ws.Cells.Item(118, 14).Value = -5114

$ws.Cells.Item(129, 8).Value = 1368
$ws.Cells.Item(129, 9).Value = 1069.25
$ws.Cells.Item(129, 10).Value = 2164.6667
$ws.Cells.Item(129, 11).Value = 3207.75
$ws.Cells.Item(129, 12).Value = 6494.000100000001
$ws.Cells.Item(129, 13).Value = 1792.25
$ws.Cells.Item(129, 14).Value = -16494.0001

$ws.Cells.Item(132, 8).Value = 2403.1646
$ws.Cells.Item(132, 9).Value = 1957.4595
$ws.Cells.Item(132, 10).Value = 8999.6
$ws.Cells.Item(132, 11).Value = 5872.3785
$ws.Cells.Item(132, 12).Value = 26998.8
$ws.Cells.Item(132, 13).Value = -3342.3785
$ws.Cells.Item(132, 14).Value = -32058.8

$ws.Cells.Item(141, 8).Value = 3967.5938
$ws.Cells.Item(141, 9).Value = 3433.2173
$ws.Cells.Item(141, 10).Value = 5333.222
$ws.Cells.Item(141, 11).Value = 10299.6519
$ws.Cells.Item(141, 12).Value = 15999.666
$ws.Cells.Item(141, 13).Value = -5119.651899999999
$ws.Cells.Item(141, 14).Value = -26359.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2687.9124
$ws.Cells.Item(32, 9).Value = 2658.4
$ws.Cells.Item(32, 10).Value = 3499.5
$ws.Cells.Item(32, 11).Value = 2658.4
$ws.Cells.Item(32, 12).Value = 3499.5
$ws.Cells.Item(32, 13).Value = -2371.4
$ws.Cells.Item(32, 14).Value = -4073.5

$ws.Cells.Item(74, 8).Value = 25642864
$ws.Cells.Item(74, 9).Value = 37038732
$ws.Cells.Item(74, 10).Value = 2156.0833
$ws.Cells.Item(74, 11).Value = 37038732
$ws.Cells.Item(74, 12).Value = 2156.0833
$ws.Cells.Item(74, 13).Value = -37037858
$ws.Cells.Item(74, 14).Value = -3904.0833

$ws.Cells.Item(77, 8).Value = 25642864
$ws.Cells.Item(77, 9).Value = 37038732
$ws.Cells.Item(77, 10).Value = 2156.0833
$ws.Cells.Item(77, 11).Value = 185193660
$ws.Cells.Item(77, 12).Value = 10780.4165
$ws.Cells.Item(77, 13).Value = -185189292
$ws.Cells.Item(77, 14).Value = -19516.4165

$ws.Cells.Item(122, 8).Value = 2944.4443
$ws.Cells.Item(122, 9).Value = 1900.2
$ws.Cells.Item(122, 10).Value = 4249.75
$ws.Cells.Item(122, 11).Value = 5700.6
$ws.Cells.Item(122, 12).Value = 12749.25
$ws.Cells.Item(122, 13).Value = -3250.6
$ws.Cells.Item(122, 14).Value = -17649.25

$ws.Cells.Item(132, 8).Value = 4808.4
$ws.Cells.Item(132, 9).Value = 4444.3823
$ws.Cells.Item(132, 10).Value = 5933.5454
$ws.Cells.Item(132, 11).Value = 13333.1469
$ws.Cells.Item(132, 12).Value = 17800.6362
$ws.Cells.Item(132, 13).Value = -10803.1469
$ws.Cells.Item(132, 14).Value = -22860.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 39109.555
$ws.Cells.Item(82, 9).Value = 4000
$ws.Cells.Item(82, 10).Value = 56664.332
$ws.Cells.Item(82, 11).Value = 4000
$ws.Cells.Item(82, 12).Value = 56664.332
$ws.Cells.Item(82, 13).Value = -3617
$ws.Cells.Item(82, 14).Value = -57430.332

$ws.Cells.Item(85, 8).Value = 39109.555
$ws.Cells.Item(85, 9).Value = 4000
$ws.Cells.Item(85, 10).Value = 56664.332
$ws.Cells.Item(85, 11).Value = 4000
$ws.Cells.Item(85, 12).Value = 56664.332
$ws.Cells.Item(85, 13).Value = -2674
$ws.Cells.Item(85, 14).Value = -59316.332

$ws.Cells.Item(99, 8).Value = 5249.25
$ws.Cells.Item(99, 9).Value = 3999.5
$ws.Cells.Item(99, 10).Value = 6499
$ws.Cells.Item(99, 11).Value = 3999.5
$ws.Cells.Item(99, 12).Value = 6499
$ws.Cells.Item(99, 13).Value = -2501.5
$ws.Cells.Item(99, 14).Value = -9495

$ws.Cells.Item(134, 8).Value = 3369.7
$ws.Cells.Item(134, 9).Value = 2586.7097
$ws.Cells.Item(134, 10).Value = 6066.6665
$ws.Cells.Item(134, 11).Value = 7760.1291
$ws.Cells.Item(134, 12).Value = 18199.9995
$ws.Cells.Item(134, 13).Value = -5225.1291
$ws.Cells.Item(134, 14).Value = -23269.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3229.5225
$ws.Cells.Item(31, 9).Value = 1446.4445
$ws.Cells.Item(31, 10).Value = 3506.2068
$ws.Cells.Item(31, 11).Value = 1446.4445
$ws.Cells.Item(31, 12).Value = 3506.2068
$ws.Cells.Item(31, 13).Value = -1151.4445
$ws.Cells.Item(31, 14).Value = -4096.2068

$ws.Cells.Item(34, 8).Value = 3229.5225
$ws.Cells.Item(34, 9).Value = 1446.4445
$ws.Cells.Item(34, 10).Value = 3506.2068
$ws.Cells.Item(34, 11).Value = 1446.4445
$ws.Cells.Item(34, 12).Value = 3506.2068
$ws.Cells.Item(34, 13).Value = -1244.4445
$ws.Cells.Item(34, 14).Value = -3910.2068

$ws.Cells.Item(105, 8).Value = 650
$ws.Cells.Item(105, 9).Value = 650
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 650
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).ClearContents()
$ws.Cells.Item(105, 14).Value = 1097

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2619.6667
$ws.Cells.Item(34, 9).Value = 100
$ws.Cells.Item(34, 10).Value = 2799.6428
$ws.Cells.Item(34, 11).Value = 300
$ws.Cells.Item(34, 12).Value = 8398.928400000001
$ws.Cells.Item(34, 13).Value = -216
$ws.Cells.Item(34, 14).Value = -8566.928400000001

$ws.Cells.Item(55, 8).Value = 1630.8
$ws.Cells.Item(55, 9).Value = 1288.5
$ws.Cells.Item(55, 10).Value = 3000
$ws.Cells.Item(55, 11).Value = 3865.5
$ws.Cells.Item(55, 12).Value = 9000
$ws.Cells.Item(55, 13).Value = -3688.5
$ws.Cells.Item(55, 14).Value = -9354

$ws.Cells.Item(121, 8).Value = 858672.8
$ws.Cells.Item(121, 9).Value = 1133508.4
$ws.Cells.Item(121, 10).Value = 34166.332
$ws.Cells.Item(121, 11).Value = 3400525.2
$ws.Cells.Item(121, 12).Value = 102498.996
$ws.Cells.Item(121, 13).Value = -3399215.2
$ws.Cells.Item(121, 14).Value = -105118.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 397.75
$ws.Cells.Item(2, 9).Value = 367.5
$ws.Cells.Item(2, 10).Value = 417.91666
$ws.Cells.Item(2, 11).Value = 367.5
$ws.Cells.Item(2, 12).Value = 417.91666
$ws.Cells.Item(2, 13).Value = -254.5
$ws.Cells.Item(2, 14).Value = -643.91666

$ws.Cells.Item(62, 8).Value = 39507
$ws.Cells.Item(62, 9).Value = 39507
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 39507
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -38821

$ws.Cells.Item(65, 8).Value = 39507
$ws.Cells.Item(65, 9).Value = 39507
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 118521
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -115089

$ws.Cells.Item(113, 8).Value = 3344.2727
$ws.Cells.Item(113, 9).Value = 1964.5
$ws.Cells.Item(113, 10).Value = 5000
$ws.Cells.Item(113, 11).Value = 1964.5
$ws.Cells.Item(113, 12).Value = 5000
$ws.Cells.Item(113, 13).Value = 205.5
$ws.Cells.Item(113, 14).Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2192.647
$ws.Cells.Item(61, 9).Value = 1885
$ws.Cells.Item(61, 10).Value = 4500
$ws.Cells.Item(61, 11).Value = 1885
$ws.Cells.Item(61, 12).Value = 4500
$ws.Cells.Item(61, 13).Value = -1683
$ws.Cells.Item(61, 14).Value = -4904

$ws.Cells.Item(63, 8).Value = 41156.4
$ws.Cells.Item(63, 9).Value = 38994
$ws.Cells.Item(63, 10).Value = 44400
$ws.Cells.Item(63, 11).Value = 38994
$ws.Cells.Item(63, 12).Value = 44400
$ws.Cells.Item(63, 13).Value = -38245
$ws.Cells.Item(63, 14).Value = -45898

$ws.Cells.Item(66, 8).Value = 41156.4
$ws.Cells.Item(66, 9).Value = 38994
$ws.Cells.Item(66, 10).Value = 44400
$ws.Cells.Item(66, 11).Value = 116982
$ws.Cells.Item(66, 12).Value = 133200
$ws.Cells.Item(66, 13).Value = -113238
$ws.Cells.Item(66, 14).Value = -140688

$ws.Cells.Item(74, 8).Value = 43500
$ws.Cells.Item(74, 9).Value = 43500
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 43500
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -42502

$ws.Cells.Item(77, 8).Value = 43500
$ws.Cells.Item(77, 9).Value = 43500
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 130500
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -125508

$ws.Cells.Item(113, 8).Value = 2192.647
$ws.Cells.Item(113, 9).Value = 1885
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 1885
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = 285
$ws.Cells.Item(113, 14).Value = -8840

$ws.Cells.Item(122, 8).Value = 3879.5217
$ws.Cells.Item(122, 9).Value = 3669.3076
$ws.Cells.Item(122, 10).Value = 4152.8
$ws.Cells.Item(122, 11).Value = 11007.9228
$ws.Cells.Item(122, 12).Value = 12458.4
$ws.Cells.Item(122, 13).Value = -8557.9228
$ws.Cells.Item(122, 14).Value = -17358.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 17125
$ws.Cells.Item(51, 9).Value = 19500
$ws.Cells.Item(51, 10).Value = 16333.333
$ws.Cells.Item(51, 11).Value = 19500
$ws.Cells.Item(51, 12).Value = 16333.333
$ws.Cells.Item(51, 13).Value = -18990
$ws.Cells.Item(51, 14).Value = -17353.333

$ws.Cells.Item(58, 8).Value = 19995.5
$ws.Cells.Item(58, 9).Value = 19995.5
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 19995.5
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -19687.5

$ws.Cells.Item(81, 8).Value = 16671609
$ws.Cells.Item(81, 9).Value = 3038.375
$ws.Cells.Item(81, 10).Value = 50008750
$ws.Cells.Item(81, 11).Value = 6076.75
$ws.Cells.Item(81, 12).Value = 100017500
$ws.Cells.Item(81, 13).Value = -5015.75
$ws.Cells.Item(81, 14).Value = -100019622

$ws.Cells.Item(84, 8).Value = 16671609
$ws.Cells.Item(84, 9).Value = 3038.375
$ws.Cells.Item(84, 10).Value = 50008750
$ws.Cells.Item(84, 11).Value = 30383.75
$ws.Cells.Item(84, 12).Value = 500087500
$ws.Cells.Item(84, 13).Value = -25079.75
$ws.Cells.Item(84, 14).Value = -500098108
